# The deck's Slide Master theme ("Integral") is swapped for the theme
# that used to live only on the Notes Master ("Office Theme"): every
# slide should now render with the Office Theme palette.
#
# All slides share one theme (ppt/theme/theme1.xml) through the single
# Slide Master, so editing the ThemeColorScheme reached from any slide
# updates it for the whole deck. The font scheme and format scheme
# (fills/lines/effects) are identical between the two themes, so only
# the twelve theme colours need to move from the Integral values to
# the Office Theme values.
#
# PowerPoint's ColorFormat.RGB is a plain Long: R + G*256 + B*65536
# (VBA's RGB() encoding), so each target hex colour is expanded below.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2,
# 7 accent3, 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
